$d = $word.ActiveDocument

# ------------------------------------------------------------------
# 1) "var type can have value of num, text, array or image or"
#    -> "var type can have value of num, txt, arr or img or"
#    (shorten "text"->"txt" and "array"->"arr" in one phrase, and
#    "image"->"img" right after it)
# ------------------------------------------------------------------
$d.Content.Find.Execute(
    "var type can have value of num, text, array", $true, $false, $false,
    $false, $false, $true, 1, $false,
    "var type can have value of num, txt, arr", 1) | Out-Null

$d.Content.Find.Execute(
    " image o", $true, $false, $false, $false, $false, $true, 1, $false,
    " img o", 1) | Out-Null

# ------------------------------------------------------------------
# 2) "g#_text" -> "g#_txt"  (only one occurrence in the document)
# ------------------------------------------------------------------
$d.Content.Find.Execute(
    "#_text", $true, $false, $false, $false, $false, $true, 1, $false,
    "#_txt", 1) | Out-Null

# ------------------------------------------------------------------
# 3) First "material,text" -> "material,txt"
#    (the "...was made of ##material,text,iron##..." sentence)
# ------------------------------------------------------------------
$first = $d.Content
$first.Find.Execute(
    "material,text", $true, $false, $false, $false, $false, $true, 1,
    $false, "", 0) | Out-Null
$first.Text = "material,txt"

# ------------------------------------------------------------------
# 4) Second "material,text" -> "material,txt"
#    (the "If the ##material,text,iron## costs ..." sentence) and
#    relocate the "_GoBack" bookmark so it now sits between the "t"
#    and "xt" of the new "txt" (i.e. after "material,t").
# ------------------------------------------------------------------
$second = $d.Range($first.End, $d.Content.End)
$second.Find.Execute(
    "material,text", $true, $false, $false, $false, $false, $true, 1,
    $false, "", 0) | Out-Null
$secondStart = $second.Start
$second.Text = "material,txt"

if ($d.Bookmarks.Exists("_GoBack")) {
    $d.Bookmarks("_GoBack").Delete()
}
$goBackPos = $secondStart + "material,t".Length
$d.Bookmarks.Add("_GoBack", $d.Range($goBackPos, $goBackPos)) | Out-Null

# ------------------------------------------------------------------
# 5) Merge the lone "##" + " " runs right after the second
#    "material,txt,iron##" into a single "## " run: extend the "##"
#    run with a trailing space, then drop the now-redundant old
#    standalone space run that follows it.
# ------------------------------------------------------------------
$hashRun = $d.Range($second.End, $d.Content.End)
$hashRun.Find.Execute(
    ",iron##", $true, $false, $false, $false, $false, $true, 1, $false,
    "", 0) | Out-Null
$hashRun = $d.Range($hashRun.End - 2, $hashRun.End)
$hashRun.Text = "## "
$oldSpace = $d.Range($hashRun.End, $hashRun.End + 1)
$oldSpace.Delete()
